$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 and Row 5 swap their prolificid/name/realeffort payload (B,C,D,F),
# while keeping their position (A index, H rank) fixed.
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "5f2c1a97a6809c060fec8820"
$ws.Range("D4").Value = "Maggie"
$ws.Range("F4").Value = 8.040674606944371

$ws.Range("B5").Value = 10
$ws.Range("C5").Value = "60a71d27a66fac796ad4de6f"
$ws.Range("D5").Value = "Jennifer"
$ws.Range("F5").Value = 8.032360915298707

# Updated realeffort scores for the remaining rows.
$ws.Range("F2").Value = 11.36491441729315
$ws.Range("F3").Value = 10.01785415257338
$ws.Range("F6").Value = 7.393130267324382
$ws.Range("F7").Value = 6.317503956260554
$ws.Range("F8").Value = 6.011467763540303
$ws.Range("F9").Value = 5.079227112452898
$ws.Range("F10").Value = 5.054581671041178
$ws.Range("F11").Value = 4.380943723260605
$ws.Range("F12").Value = 2.380573180982641
$ws.Range("F13").Value = 0.2369653110842641
